$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Widen column A to fit the longer program names
# (target stored width 34.6640625; the runtime's character->pixel
# rounding only lands on multiples of 1/6, so 33.83 is the closest
# achievable ColumnWidth, landing on stored width 34.666666666666664)
$ws.Columns.Item(1).ColumnWidth = 33.83

# Update the Yes/No answers for the existing programs (rows 3-11)
$ws.Range("B3").Value = "No"
$ws.Range("B4").Value = "Yes"
$ws.Range("B5").Value = "Yes"
$ws.Range("B6").Value = "Yes"
$ws.Range("B7").Value = "Yes"
$ws.Range("B8").Value = "Yes"
$ws.Range("B9").Value = "No"
$ws.Range("B10").Value = "No"
$ws.Range("B11").Value = "No"

# Add the new programs
$ws.Range("A12").Value = "KIT_ME"
$ws.Range("B12").Value = "Yes"

$ws.Range("A13").Value = "TU_DORTMUND_ROBOTICS"
$ws.Range("B13").Value = "Yes"

$ws.Range("A14").Value = "RWTH_Aachen_ROBOTICS"
$ws.Range("B14").Value = "Yes"

# Extend the data validation list (Yes/No) to cover the new rows
$newValidationRange = $ws.Range("B1:B14")
$newValidationRange.Validation.Delete()
$newValidationRange.Validation.Add(3, 1, 1, '"Yes,No"')
$newValidationRange.Validation.IgnoreBlank = $true
$newValidationRange.Validation.InCellDropdown = $true
$newValidationRange.Validation.ShowInput = $true
$newValidationRange.Validation.ShowError = $true

# Move the selection like in the saved file
[void]$ws.Range("E5").Select()
